$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Cases" query text (was previously in column A/B before the new TabName
# column was inserted at the front).
$casesQuery = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)`n    WHERE c.race = " + [char]34 + "WHITE" + [char]34 + "`nWITH DISTINCT c, a, ct`nRETURN `n    COALESCE(c.case_id, '') AS ``Case ID``,`n    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n    COALESCE(a.arm_id, '') AS ``Arm``,`n    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,`n    COALESCE(c.disease, '') AS ``Diagnosis``,`n    COALESCE(c.gender, '') AS ``Gender``,`n    COALESCE(c.race, '') AS ``Race``,`n    COALESCE(c.ethnicity, '') AS ``Ethnicity``"

# New "StatQuery" text.
$statQuery = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n    WHERE c.race = " + [char]34 + "WHITE" + [char]34 + "`nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`nRETURN `n    COUNT(DISTINCT f) AS number_of_files,`n    COUNT(DISTINCT c.case_id) AS number_of_cases,`n    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

# Insert a new "TabName" column at the very front, shifting the existing
# four columns (query, StatQuery, dbExcel, WebExcel) one to the right.
# NOTE: inserting like this preserves the exact original column widths of
# the shifted columns (B:E) untouched - only the brand-new column A needs an
# explicit width below.
$ws.Columns("A").Insert()

# Row 1 - headers
$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# Row 2 - values
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = "TC06_Trials_Filter_Race-White_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC06_Trials_Filter_Race-White_WebData.xlsx"

# Column width for the newly-inserted "TabName" column only; the other
# columns (B:E) already carry their correct original widths forward from
# the insert above, so we deliberately leave them alone.
$ws.Columns("A").ColumnWidth = 8

# Row height for the data row (wrapped long query text).
$ws.Rows("2").RowHeight = 174

# Wrap text + style for the two long query cells.
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true

# Restore the selection to match the saved workbook state.
$ws.Range("B10").Select() | Out-Null
